$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = "FreeRTOSv202112.00"
$ws.Range("C9").Value = "FreeRTOSv202112.00 源码"

$ws.Range("G11").Select()
